$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# Insert a new row at 178, pushing existing rows 178+ down to 179+
$ws.Rows.Item(178).Insert()

# Row 165 already carries the exact cell-style combination this new row needs
# (fills/borders for the 55900xxx band), so copy formats only from it.
$ws.Range("A165:Z165").Copy()
$ws.Range("A178:Z178").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 178 with the "chenmo" (Silent) monster skill data
$ws.Range("A178").Value = 55900046
$ws.Range("B178").Value = "沉默"
$ws.Range("C178").Value = "特效"
$ws.Range("D178").Value = "NAR"
$ws.Range("E178").Value = 20
$ws.Range("H178").Value = "foreach(IMonster mon in s.Map.GetRangeMonster(s.IsLeft,sp.Target,sp.Shape,sp.Range,s.Position)) if(mon.Id!=s.Id) mon.Silent();"
$ws.Range("Q178").Value = "Active"
$ws.Range("R178").Formula = '="true"'
$ws.Range("R178").Copy()
$ws.Range("R178").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("S178").Value = "召唤时沉默范围内所有目标"
$ws.Range("U178").Value = "yellowflash"
$ws.Range("V178").Value = "yellowflash"
$ws.Range("X178").Value = 25
$ws.Range("Y178").Value = "chenmo"

$ws.Rows.Item(178).RowHeight = 72

# Update the table range to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:Z194"))

# Update dimension-related view state
$ws.Application.ActiveWindow.ScrollRow = 176
$ws.Range("E178").Select()
